# Weekly update: prepend a new week's Hass/Edranol avocado price rows
# (Terminal La Palmera de La Serena) ahead of the existing history.
# This pushes all existing data rows (1068:1085) down by 5, to (1073:1090).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 fresh rows above the current first data block of this commodity
# (row 1068), shifting all subsequent rows down by 5.
$ws.Rows("1068:1072").Insert()

# New week's rows (date serial 44890 = 2022-11-25), same market/region/
# product metadata as the rest of the sheet.
$newRows = @(
    @{ Row=1068; K="Edranol"; L="Primera"; M=200; N=1800;  O=1900;  P=1850;  Q="$/kilo (en caja de 17 kilos)"; S=1850; T=1 },
    @{ Row=1069; K="Edranol"; L="Segunda"; M=160; N=1600;  O=1700;  P=1650;  Q="$/kilo (en caja de 17 kilos)"; S=1650; T=1 },
    @{ Row=1070; K="Hass";    L="Primera"; M=240; N=2500;  O=2600;  P=2550;  Q="$/kilo (en caja de 17 kilos)"; S=2550; T=1 },
    @{ Row=1071; K="Hass";    L="Segunda"; M=200; N=2200;  O=2300;  P=2250;  Q="$/kilo (en caja de 17 kilos)"; S=2250; T=1 },
    @{ Row=1072; K="Hass";    L="Tercera"; M=160; N=1800;  O=1900;  P=1850;  Q="$/kilo (en caja de 17 kilos)"; S=1850; T=1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 8
    $ws.Cells.Item($row, 2).Value  = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value  = "Coquimbo"
    $ws.Cells.Item($row, 4).Value  = 44890
    $ws.Cells.Item($row, 5).Value  = 4
    $ws.Cells.Item($row, 6).Value  = "Fruta"
    $ws.Cells.Item($row, 7).Value  = 100106
    $ws.Cells.Item($row, 8).Value  = "Oleaginosos"
    $ws.Cells.Item($row, 9).Value  = 100106002
    $ws.Cells.Item($row, 10).Value = "Palta"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
